$wb = $excel.ActiveWorkbook

# --- "Variables" sheet: merge "collection event.resource" + "collection event.name"
# columns into a single "collection event" column. The resource half (col I) is
# redundant (always equals the row's own "resource" value), so the edit simply
# drops that column, leaving the previous "collection event.name" column (now
# shifted left into column I) and relabels its header.
$wsVariables = $wb.Worksheets.Item("Variables")
[void]$wsVariables.Columns.Item(9).EntireColumn.Delete()
$wsVariables.Range("I1").Value = "collection event"

# --- "Repeated variables" sheet: same merge, columns G/H instead of I/J.
$wsRepeated = $wb.Worksheets.Item("Repeated variables")
[void]$wsRepeated.Columns.Item(7).EntireColumn.Delete()
$wsRepeated.Range("G1").Value = "collection event"

# --- Update the saved cursor/selection on every sheet, finishing on
# "Datasets" so it ends up the active tab.
[void]$wsVariables.Range("A2").Select()

$wsValues = $wb.Worksheets.Item("Variable values")
[void]$wsValues.Range("A2").Select()

[void]$wsRepeated.Range("A2").Select()

$wsDatasets = $wb.Worksheets.Item("Datasets")
[void]$wsDatasets.Range("B2").Select()
